$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features = @(
    'OverallQual',
    'GarageCars',
    'FullBath',
    'GarageType_Detchd',
    'MSZoning_RM',
    'CentralAir',
    'TotalBsmtSF',
    'Fireplaces',
    'Neighborhood_Crawfor',
    '2ndFlrSF',
    'GarageType_Attchd',
    'BsmtFinSF1',
    'BuildingAge',
    'HeatingQC',
    'BsmtExposure_Gd',
    'MSZoning_RL',
    'SaleType_New',
    'LotArea',
    'Neighborhood_OldTown',
    'GarageFinish_Unf',
    'Exterior1st_BrkFace',
    'BsmtFinType1_GLQ',
    'GarageType_No Garage',
    'MasVnrArea',
    'OverallCond',
    'GarageQual',
    'ScreenPorch',
    'HouseStyle_1Story',
    'SaleCondition_Partial',
    'HouseStyle_2Story',
    'Foundation_PConc',
    'RoofMatl_Tar&Grv',
    'MSZoning_FV',
    'LotFrontage',
    'OpenPorchSF',
    'Neighborhood_NoRidge',
    'LotConfig_CulDSac',
    'RoofStyle_Gambrel',
    'Neighborhood_Timber',
    'Neighborhood_StoneBr',
    'SaleCondition_Normal',
    'BsmtCond',
    'BedroomAbvGr',
    'KitchenAbvGr',
    'SaleType_WD',
    'Functional_Typ',
    'Exterior1st_CemntBd',
    'Neighborhood_BrkSide',
    'LotShape_Reg',
    'WoodDeckSF',
    'Condition1_Norm',
    'Exterior2nd_VinylSd',
    'SaleCondition_Family',
    'Neighborhood_Somerst',
    'RoofStyle_Gable',
    'LandContour_Low',
    'Neighborhood_NridgHt',
    'BsmtFinSF2',
    'Exterior1st_WdShing',
    'BsmtHalfBath'
)

for ($i = 0; $i -lt $features.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $features[$i]
}
